$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now represents the "جيانا" sample instead of "ديما".
# Updating A2 automatically: recalculates the B2 formula's cached value,
# reindexes/garbage-collects the shared-string table (dropping "ديما"),
# and leaves D2's literal value (1) untouched.
$ws.Range("A2").Value = "جيانا"

# Row 3 (previously the "جيانا" row) becomes a blank templated row,
# matching the look of the old blank row 4 - just clear its contents,
# keep the existing cell formatting/style.
$ws.Range("A3:D3").ClearContents()

# Duplicate that blank templated row (row 3) into a new row 5, pushing
# the filler rows below down by one (this is what creates the new
# trailing filler row 18 and bumps the sheet dimension to D18).
$ws.Rows(3).Copy()
$ws.Rows(4).Insert()

# The old blank templated row 4 is no longer needed (its job is now done
# by rows 3 and 5) - clear it completely so it collapses back out of the
# worksheet entirely.
$ws.Range("A4:D4").Clear()

# Reflect the new selection left by the edit (a full selector row).
$ws.Range("A3:E3").Select()
